# Update the "2024" worksheet:
#   - A new September transaction ("cashback tap pay visa" on 2024-09-19 15:37:45)
#     is logged at the top of the September list (columns R:S), pushing all
#     later September entries down by one row.
#   - A new (empty) category row is inserted above the "Braodband" row in
#     column A, pushing every subsequent category label down by one row and
#     extending the sheet by a row (the last label, "Broadband", ends up on
#     the new row 172).
#
# Because this runtime's Range.Insert() shifts an entire row regardless of
# which columns are referenced, the shift is instead performed manually by
# reading the affected column(s) into an array, shifting the values down by
# one slot, and writing the new values (including the newly logged
# transaction) back out in one shot.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# --- Column A ("Group" categories): shift rows 7-171 down to rows 8-172,
#     leaving row 7 blank. ---
$colA  = $ws.Range("A7:A171").Value2
$countA = $colA.GetLength(0)
$sizeA  = $countA + 1
$newA = New-Object 'object[,]' $sizeA,1
$newA[0,0] = $null
for ($i = 1; $i -le $countA; $i++) {
    $newA[$i,0] = $colA[$i,1]
}
$ws.Range("A7:A172").Value2 = $newA

# --- Columns R:S (September_Details / September_Date): shift rows 3-162
#     down to rows 4-163, then record the newly logged transaction on row 3. ---
$colRS  = $ws.Range("R3:S162").Value2
$countRS = $colRS.GetLength(0)
$sizeRS  = $countRS + 1
$newRS = New-Object 'object[,]' $sizeRS,2
$newRS[0,0] = "cashback tap pay visa"
$newRS[0,1] = "2024-09-19 15:37:45"
for ($i = 1; $i -le $countRS; $i++) {
    $newRS[$i,0] = $colRS[$i,1]
    $newRS[$i,1] = $colRS[$i,2]
}
$ws.Range("R3:S163").Value2 = $newRS
